$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'62.517.59"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = "'3.434.20"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'406.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'131.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = "'0.693"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('E10').Value = '  +7.52%  '
$ws.Range('D11').Value = "'42.11"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'8.41"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'19.86"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = "'3.428.72"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').Value = "'62.414.47"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = "'11.59"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = "'1.02"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').Value = "'0.0000150"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.75%  '
$ws.Range('D20').Value = "'3.18"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').Value = "'84.57"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('D22').Value = "'312.24"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').Value = "'12.84"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('D24').Value = "'3.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('D26').Value = "'29.72"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('D27').Value = "'8.11"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.39%  '
$ws.Range('D28').Value = "'7.76"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('E29').Value = '  +4.62%  '
$ws.Range('D30').Value = "'44.89"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.63%  '
$ws.Range('D31').Value = "'0.173"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('D33').Value = "'11.38"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'0.0484"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('D36').Value = "'51.77"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').Value = "'0.998"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('D39').Value = "'0.321"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.78%  '
$ws.Range('D40').Value = "'3.29"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.40%  '
$ws.Range('D41').Value = "'142.63"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.69%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').Value = "'1.98"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('D44').Value = "'3.93"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').Value = "'16.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').Value = "'21.36"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('D48').Value = "'2.106.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('B49').Value = 'OceanProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D49').Value = "'1.12"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +31.54%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = "'1.97"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.52%  '
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').Value = "'2.31"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.65%  '
